# "added 4wk low sales check"
# Update the forecast comparison figures (MyForecast, Inventory Coverage,
# Seasonality Index) and the derived Summary totals/min that change as a
# result of adding the new 4-week-low-sales check.

$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" --------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("L2").Value = 1.18

$wsForecast.Range("D3").Value = 13
$wsForecast.Range("H3").Value = 6.54
$wsForecast.Range("L3").Value = 0.8

$wsForecast.Range("D4").Value = 13
$wsForecast.Range("H4").Value = 5.54
$wsForecast.Range("L4").Value = 1.15

$wsForecast.Range("H5").Value = 4.54
$wsForecast.Range("L5").Value = 0.84

$wsForecast.Range("D6").Value = 12
$wsForecast.Range("H6").Value = 3.83
$wsForecast.Range("L6").Value = 0.95

$wsForecast.Range("D7").Value = 12
$wsForecast.Range("H7").Value = 2.83
$wsForecast.Range("L7").Value = 1

$wsForecast.Range("H8").Value = 1.69
$wsForecast.Range("L8").Value = 1.02

$wsForecast.Range("H9").Value = 0.6899999999999999
$wsForecast.Range("L9").Value = 0.87

$wsForecast.Range("D10").Value = 13
$wsForecast.Range("L10").Value = 1.04

$wsForecast.Range("D11").Value = 12
$wsForecast.Range("L11").Value = 1.07

$wsForecast.Range("D12").Value = 12
$wsForecast.Range("L12").Value = 1.12

$wsForecast.Range("D13").Value = 12
$wsForecast.Range("L13").Value = 0.84

$wsForecast.Range("D14").Value = 12
$wsForecast.Range("L14").Value = 0.96

$wsForecast.Range("D15").Value = 12
$wsForecast.Range("L15").Value = 1.12

$wsForecast.Range("D16").Value = 12
$wsForecast.Range("L16").Value = 0.82

$wsForecast.Range("D17").Value = 12
$wsForecast.Range("L17").Value = 0.8100000000000001

# --- Sheet "Summary" ----------------------------------------------------
# Values on this sheet are stored as text (e.g. "191"), so a leading
# apostrophe is used to force a text entry instead of a number. The
# formatting picked up from the text-entry quote-prefix is cleared right
# away so the cell style stays identical to the original (General/no style).
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").Value = "'200"
$wsSummary.Range("B9").ClearFormats()

$wsSummary.Range("B10").Value = "'103"
$wsSummary.Range("B10").ClearFormats()

$wsSummary.Range("B11").Value = "'53"
$wsSummary.Range("B11").ClearFormats()

$wsSummary.Range("B14").Value = "'12"
$wsSummary.Range("B14").ClearFormats()
